$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "PEDC" column is inserted before the existing "TOTAL_BCQ" column:
#   column E, formerly TOTAL_BCQ, now holds the new PEDC values
#   column F is the (new) TOTAL_BCQ column = SCPC + KSPC + EDC + PEDC
$ws.Range("E1").Value = "PEDC"
$ws.Range("F1").Value = "TOTAL_BCQ"

# Per-row values for columns B..F (SCPC, KSPC, EDC, PEDC, TOTAL_BCQ), indexed by HOUR (1-24, sheet rows 2-25)
$rows = @(
    @(12500, 10000, 20000, 5000, 47500),
    @(12500, 10000, 10000, 5000, 37500),
    @(12500, 10000, 0, 5000, 27500),
    @(12500, 10000, 0, 5000, 27500),
    @(12500, 10000, 0, 5000, 27500),
    @(12500, 10000, 0, 5000, 27500),
    @(12500, 10000, 0, 5000, 27500),
    @(12500, 10000, 0, 5000, 27500),
    @(25000, 20000, 20000, 5000, 70000),
    @(25000, 20000, 20000, 5000, 70000),
    @(25000, 20000, 20000, 10000, 75000),
    @(25000, 20000, 20000, 10000, 75000),
    @(25000, 20000, 20000, 10000, 75000),
    @(25000, 20000, 20000, 10000, 75000),
    @(25000, 20000, 20000, 10000, 75000),
    @(25000, 20000, 20000, 10000, 75000),
    @(25000, 20000, 20000, 10000, 75000),
    @(25000, 20000, 20000, 10000, 75000),
    @(25000, 20000, 20000, 10000, 75000),
    @(25000, 20000, 20000, 10000, 75000),
    @(25000, 20000, 20000, 10000, 75000),
    @(25000, 20000, 20000, 10000, 75000),
    @(25000, 20000, 20000, 10000, 75000),
    @(12500, 10000, 20000, 5000, 47500)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $vals = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
}
